# Energy InNovation ESF workbook restructuring
# 1. Duplicate the "ESF" sheet -> "ESF (2)" (placed immediately after "ESF")
# 2. Wrap the existing "ESF" formulas in MAX(0, ...)
# 3. Rewrite the new "ESF (2)" formulas to reference 'AEO 2021 Table 41' row 26
#    divided by its own $AJ$26 absolute total
# 4. Insert a new blank worksheet named "Sheet2" immediately before "ESF"
# 5. Restore/refresh a couple of cosmetic sheet-view selections
#
# NOTE: worksheet object references returned by $wb.Worksheets.Item(...) can
# go stale (start pointing at the wrong sheet) once the sheet collection is
# reshuffled by an Add/Copy/Delete/Move, so we always re-fetch sheets by
# name right before we touch them.

$wb = $excel.ActiveWorkbook

# --- Step 1: duplicate ESF -> "ESF (2)" placed right after ESF ---------
$esf = $wb.Worksheets.Item("ESF")
$esf.Copy([System.Reflection.Missing]::Value, $esf)

# --- Step 2: wrap ESF's row-2 formulas in MAX(0, ...) -------------------
$esfCols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF")

$esf = $wb.Worksheets.Item("ESF")
foreach ($col in $esfCols) {
    $cell = $esf.Range($col + "2")
    $oldFormula = $cell.Formula
    $inner = $oldFormula.Substring(1)
    $cell.Formula = "=MAX(0," + $inner + ")"
}

# --- Step 3: rewrite ESF (2)'s row-2 formulas ----------------------------
$t41cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ")

$esf2 = $wb.Worksheets.Item("ESF (2)")
for ($i = 0; $i -lt $esfCols.Length; $i++) {
    $destCol = $esfCols[$i]
    $srcCol = $t41cols[$i]
    $cell = $esf2.Range($destCol + "2")
    $cell.Formula = "=(1-'AEO 2021 Table 41'!" + $srcCol + "26/'AEO 2021 Table 41'!`$AJ`$26)"
}

# --- Step 4: new blank sheet "Sheet2" inserted right before ESF --------
$esf = $wb.Worksheets.Item("ESF")
$blank = $wb.Worksheets.Add($esf)
$blank.Name = "Sheet2"
$blank.Range("I33").Select()

# --- Step 5: cosmetic selection fixes ------------------------------------
$t41 = $wb.Worksheets.Item("AEO 2021 Table 41")
$t41.Range("G37").Select()

$esf2 = $wb.Worksheets.Item("ESF (2)")
$esf2.Range("I4").Select()

$esf = $wb.Worksheets.Item("ESF")
$esf.Activate()
$esf.Range("B2:AF2").Select()
